$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out counter values that no longer apply after cleaning up the
# data processing (operation counters reset to 0).
$ws.Range("B2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
